# Apply the "Modified first 10 rows based on recent estimates" edit.
#
# Changes:
#  1. CPU time per event (row 9) now derived from Oct. 2016 David benchmark
#     rates (250 Hz / 16 threads Ivy Bridge, 340 Hz / 24 threads Haswell,
#     600 Hz / 36 threads Broadwell) instead of a flat 1/22 s/event.
#  2. Raw event size (row 11) is now a formula built from the Spring 2016
#     data-driven estimate (16.1/11.5 kB fixed + 4.6/2.3 kB per 10^7 g/s,
#     scaled per-column) instead of a flat 18000-byte literal.
#  3. Updated the accompanying comment text in column F for both rows.
#  4. Cosmetic: window x-position and the active selection on the "model"
#     sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# --- Row 9: CPU time per event -------------------------------------------
$ws.Range("B9").Formula = "=1/(340/24)"
$ws.Range("C9").Formula = "=1/(340/24)"
$ws.Range("D9").Formula = "=1/(600/36)"
$ws.Range("F9").Value = "Oct. 2016 David benchmark gave 250Hz for 16 threads on Ivy Bridge, 340Hz for 24 threads on Haswell, 600Hz for 36 threads on Broadwell"

# --- Row 11: raw event size ----------------------------------------------
$ws.Range("B11").Formula = "=(16.1+4.6*1)*1000"
$ws.Range("C11").Formula = "=(11.5+2.3*2)*1000"
$ws.Range("D11").Formula = "=(11.5+2.3*5)*1000"
$ws.Range("F11").Value = "size of a single raw event. Actual Spring 2016 data is 16.4kB+4.6kB/10^7 g/s. Estimate from Spring 2016 data for reduced windows is 11.5kB + 0.23kB/10^7 g/s"

# --- Cosmetic: window position + active selection -------------------------
# (window screen x-position isn't part of the writable Excel object model
# surface here, so only the cell-selection change is applied)
$excel.ActiveWindow.Left = 600

$ws.Activate()
$ws.Range("A15").Select()

$wb.Application.Calculate()
